# Add a new "BOUNDARY" adversarial-attack result block (res boundary attack,
# SEED 42) to the results sheet: an 8-column section (eps 0.01..0.20) in
# columns AK:AR mirroring the existing NA/PAST/REV/REV_BIM/FGSM_SURRO/FGSM
# blocks, with MAE/RMSE/SIM rows for each of the LSTM/RNN/GRU models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: BOUNDARY attack header (AK1:AR1), merged, styled like the AC1:AJ1 block
$ws.Range("AK1").Value = "BOUNDARY"

# Row 2: epsilon labels for the BOUNDARY block (text, not numbers)
$ws.Range("AK2").Value = "'0.01"
$ws.Range("AL2").Value = "'0.02"
$ws.Range("AM2").Value = "'0.03"
$ws.Range("AN2").Value = "'0.04"
$ws.Range("AO2").Value = "'0.05"
$ws.Range("AP2").Value = "'0.07"
$ws.Range("AQ2").Value = "'0.10"
$ws.Range("AR2").Value = "'0.20"

# Data rows for the BOUNDARY block: MAE (4,7,10), RMSE (5,8,11), SIM (6,9,12)
$ws.Range("AK4").Value = 429.9866724268595
$ws.Range("AL4").Value = 432.1328560256958
$ws.Range("AM4").Value = 438.2829358355204
$ws.Range("AN4").Value = 441.9419984690348
$ws.Range("AO4").Value = 458.2547001711528
$ws.Range("AP4").Value = 467.7945947647095
$ws.Range("AQ4").Value = 476.5419515482585
$ws.Range("AR4").Value = 652.0824718093872

$ws.Range("AK5").Value = 526.879090701764
$ws.Range("AL5").Value = 528.7592190548929
$ws.Range("AM5").Value = 532.465842137164
$ws.Range("AN5").Value = 539.6932150139494
$ws.Range("AO5").Value = 554.7810853149108
$ws.Range("AP5").Value = 565.2909471172483
$ws.Range("AQ5").Value = 569.3211817184498
$ws.Range("AR5").Value = 824.4355533053063

$ws.Range("AK6").Value = 0.9989798422214393
$ws.Range("AL6").Value = 0.998972084297574
$ws.Range("AM6").Value = 0.998956535769256
$ws.Range("AN6").Value = 0.9989305246276561
$ws.Range("AO6").Value = 0.9988709320757765
$ws.Range("AP6").Value = 0.9988615007234319
$ws.Range("AQ6").Value = 0.9988037831648963
$ws.Range("AR6").Value = 0.9974962230004524

$ws.Range("AK7").Value = 397.4771784273784
$ws.Range("AL7").Value = 399.6712321281433
$ws.Range("AM7").Value = 400.1091247812907
$ws.Range("AN7").Value = 412.6248471069336
$ws.Range("AO7").Value = 419.6750252087911
$ws.Range("AP7").Value = 464.7955890782674
$ws.Range("AQ7").Value = 492.2165950520833
$ws.Range("AR7").Value = 778.9044143040975

$ws.Range("AK8").Value = 516.8373745255692
$ws.Range("AL8").Value = 518.6026726811065
$ws.Range("AM8").Value = 519.8309088780161
$ws.Range("AN8").Value = 528.8819728631149
$ws.Range("AO8").Value = 536.9590687835486
$ws.Range("AP8").Value = 594.5860379004587
$ws.Range("AQ8").Value = 651.2489338147866
$ws.Range("AR8").Value = 959.3623680262759

$ws.Range("AK9").Value = 0.9992333950942115
$ws.Range("AL9").Value = 0.9992350673850453
$ws.Range("AM9").Value = 0.9992230233412658
$ws.Range("AN9").Value = 0.9992157602768351
$ws.Range("AO9").Value = 0.9991527386700584
$ws.Range("AP9").Value = 0.9988674522484662
$ws.Range("AQ9").Value = 0.998639055191368
$ws.Range("AR9").Value = 0.9969367382912717

$ws.Range("AK10").Value = 306.5038069725036
$ws.Range("AL10").Value = 307.0465618006388
$ws.Range("AM10").Value = 312.7120273844401
$ws.Range("AN10").Value = 340.8137208811442
$ws.Range("AO10").Value = 338.9568512598673
$ws.Range("AP10").Value = 381.1121523920695
$ws.Range("AQ10").Value = 394.5817514101664
$ws.Range("AR10").Value = 706.387831103007

$ws.Range("AK11").Value = 419.5447428040518
$ws.Range("AL11").Value = 414.9333354934766
$ws.Range("AM11").Value = 423.5129873958635
$ws.Range("AN11").Value = 460.5121916836027
$ws.Range("AO11").Value = 454.1124385065061
$ws.Range("AP11").Value = 514.3530651598844
$ws.Range("AQ11").Value = 536.4903620755837
$ws.Range("AR11").Value = 891.6669043707768

$ws.Range("AK12").Value = 0.9993254988044271
$ws.Range("AL12").Value = 0.999341280220878
$ws.Range("AM12").Value = 0.9993120574863638
$ws.Range("AN12").Value = 0.999181721992064
$ws.Range("AO12").Value = 0.9992135930996382
$ws.Range("AP12").Value = 0.998969762262799
$ws.Range("AQ12").Value = 0.9988815050240373
$ws.Range("AR12").Value = 0.9968422392808229

# Merge the new header cell the same way the other attack headers are merged
$ws.Range("AK1:AR1").Merge() | Out-Null

# Copy cell formatting (border/bold/alignment) from the neighboring FGSM header/epsilon
# row cells onto the new BOUNDARY header row so the new block matches the others.
$ws.Range("AJ1").Copy() | Out-Null
$ws.Range("AK1:AR1").PasteSpecial(-4122) | Out-Null
$ws.Range("AJ2").Copy() | Out-Null
$ws.Range("AK2:AR2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Tiny last-ULP recompute deltas in the existing SIM rows caused by the new block
$ws.Range("E6").Value = 0.9989914276568976
$ws.Range("H6").Value = 0.9989185339197603
$ws.Range("J6").Value = 0.9987771386465392
$ws.Range("K6").Value = 0.9984846802261947
$ws.Range("L6").Value = 0.9981562031166326
$ws.Range("M6").Value = 0.9989853082799081
$ws.Range("R6").Value = 0.9989049101038218
$ws.Range("V6").Value = 0.9982819839665024
$ws.Range("W6").Value = 0.9977865002336597
$ws.Range("AC6").Value = 0.998613994532612
$ws.Range("AF6").Value = 0.9968637635670509
$ws.Range("AG6").Value = 0.9960858956958423
$ws.Range("C9").Value = 0.9992496037173066
$ws.Range("D9").Value = 0.9992423920092128
$ws.Range("F9").Value = 0.999157243917792
$ws.Range("G9").Value = 0.9990948497268954
$ws.Range("I9").Value = 0.9989387630646054
$ws.Range("K9").Value = 0.9983435919244245
$ws.Range("P9").Value = 0.9990567408627172
$ws.Range("Q9").Value = 0.9990009800176723
$ws.Range("W9").Value = 0.9986066879513898
$ws.Range("Y9").Value = 0.9979590840130074
$ws.Range("AF9").Value = 0.9983719084275864
$ws.Range("AI9").Value = 0.9954034259128793
$ws.Range("AJ9").Value = 0.9852036529759525
$ws.Range("F12").Value = 0.9992553124273974
$ws.Range("G12").Value = 0.9992451854697948
$ws.Range("H12").Value = 0.9992261846084264
$ws.Range("J12").Value = 0.9990903949289373
$ws.Range("M12").Value = 0.999352529974611
$ws.Range("R12").Value = 0.9991805740570844
$ws.Range("S12").Value = 0.9989760492569433
$ws.Range("W12").Value = 0.9986672514997205
$ws.Range("Z12").Value = 0.9968242447811232
$ws.Range("AB12").Value = 0.9865924354164557
$ws.Range("AI12").Value = 0.9930027546959711
